# Insert a new data row at row 142 (shifting the existing rows 142-237
# down to 143-238), mirroring a new weekly price observation for
# Jengibre @ Vega Modelo de Temuco being added to the top of the
# historical series kept in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 142:237 down one row to make room for the new record.
$ws.Rows.Item(142).Insert()

# Populate the newly inserted row 142 with the new observation.
$ws.Range("A142").Value = 10
$ws.Range("B142").Value = "Vega Modelo de Temuco"
$ws.Range("C142").Value = "La Araucanía"
$ws.Range("D142").Value = 44942
$ws.Range("E142").Value = 9
$ws.Range("F142").Value = 100114007
$ws.Range("G142").Value = "Jengibre"
$ws.Range("H142").Value = "Sin especificar"
$ws.Range("I142").Value = "Primera"
$ws.Range("J142").Value = 40
$ws.Range("K142").Value = 20000
$ws.Range("L142").Value = 20000
$ws.Range("M142").Value = 20000
$ws.Range("N142").Value = "$/caja 13 kilos"
$ws.Range("O142").Value = "Perú"
$ws.Range("P142").Value = 1538
$ws.Range("Q142").Value = 13
$ws.Range("R142").Value = "Hortaliza"
